$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1: new header cell, same style as A1/B1 (bold, centered, bordered)
$ws.Range("C1").Value = "Share `$\le`$ 1 Rating From Actual"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# C2:C4: numeric-looking values stored as text
$ws.Range("C2:C4").NumberFormat = "@"
$ws.Range("C2").Value = "0.7687"
$ws.Range("C3").Value = "0.9254"
$ws.Range("C4").Value = "0.9353"

# C5: present but empty, to extend the used range to C5
$ws.Cells.Item(5, 3).Borders.LineStyle = 0
